$wb = $excel.ActiveWorkbook

# --- NamedThing: drop the trailing type_value column (D1) ---
$ws = $wb.Worksheets.Item("NamedThing")
$ws.Range("D1").ClearContents()

# --- Observation -> Person: becomes a simple person record ---
$ws = $wb.Worksheets.Item("Observation")
$ws.Name = "Person"
$ws.Range("D1").ClearContents()
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "first_name"
$ws.Range("C1").Value = "last_name"

# --- NonProcess -> PersonCollection: becomes a homogeneous collection of people ---
$ws = $wb.Worksheets.Item("NonProcess")
$ws.Name = "PersonCollection"
$ws.Range("B1").ClearContents()
$ws.Range("C1").ClearContents()
$ws.Range("D1").ClearContents()
$ws.Range("A1").Value = "people"

# --- MaterialEntity: drop the trailing type_value column (F1) ---
$ws = $wb.Worksheets.Item("MaterialEntity")
$ws.Range("F1").ClearContents()

# --- SoilSample: drop the trailing type_value column (F1) ---
$ws = $wb.Worksheets.Item("SoilSample")
$ws.Range("F1").ClearContents()

# --- DnaExtract: drop the trailing type_value column (F1) ---
$ws = $wb.Worksheets.Item("DnaExtract")
$ws.Range("F1").ClearContents()

# --- InformationArtifact: drop the trailing type_value column (G1) ---
$ws = $wb.Worksheets.Item("InformationArtifact")
$ws.Range("G1").ClearContents()

# --- Process: drop the trailing type_value column (F1) ---
$ws = $wb.Worksheets.Item("Process")
$ws.Range("F1").ClearContents()
